# Generate Report for Handoff
# A new file (77a4a6e8-9225-47f5-b53f-bfe9c910a317) has reached "Ready for
# handoff" status. Insert it as the newest row (row 2) on every sheet,
# pushing the previous newest entry (96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8)
# down to row 3.

$wb = $excel.ActiveWorkbook

$oldGuid = "96e8afa5-9aa7-4dfe-8212-60b1e5cc62e8"
$newGuid = "77a4a6e8-9225-47f5-b53f-bfe9c910a317"
$oldHash = "0012e40d796e5c6f54b3c87d5af7bf616b8ae37b"
$newHash = "61be0d3a06f1e2c0b0a43f3cebec17b440913912"

$mdBlobSha = "3073ab0b57de045a493efbd6ee1c3f4aaaaaa857"
$zhHandoffSha = "a8d307b6cb8b29c0798d457611105d1993f1f720"
$deHandoffSha = "335346405c8c72e7301abad97d5697d2151ff791"

# ============================================================
# Sheet 1: Overview
# ============================================================
$ws1 = $wb.Worksheets.Item(1)
$ws1.Hyperlinks.Delete()

# Push the existing data row down to make room for the new entry.
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = "$oldGuid.md"
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"
$ws1.Range("D3").Value = "2016-03-22 00:33:37"

$ws1.Range("A2").Value = "$newGuid.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"
$ws1.Range("D2").Value = "2016-03-22 00:34:05"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$oldGuid.md", "", "", "$oldGuid.md")

# ============================================================
# Sheet 2: zh-cn
# ============================================================
$ws2 = $wb.Worksheets.Item(2)
$ws2.Hyperlinks.Delete()

$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = "$oldGuid.md"
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "$oldGuid.$oldHash.zh-cn.xlf"
$ws2.Range("E3").Value = "2016-03-22 00:33:33"
$ws2.Range("H3").Value = "0001-01-01 00:00:00"
$ws2.Range("J3").Value = "Include"

$ws2.Range("A2").Value = "$newGuid.md"
$ws2.Range("B2").Value = ".md"
$ws2.Range("C2").Value = "Ready for handoff"
$ws2.Range("D2").Value = "$newGuid.$newHash.zh-cn.xlf"
$ws2.Range("E2").Value = "2016-03-22 00:34:00"
$ws2.Range("H2").Value = "0001-01-01 00:00:00"
$ws2.Range("J2").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$newGuid.$newHash.zh-cn.xlf", "", "", "$newGuid.$newHash.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$oldGuid.md", "", "", "$oldGuid.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldGuid.$oldHash.zh-cn.xlf", "", "", "$oldGuid.$oldHash.zh-cn.xlf")

# ============================================================
# Sheet 3: de-de
# ============================================================
$ws3 = $wb.Worksheets.Item(3)
$ws3.Hyperlinks.Delete()

$ws3.Rows.Item(3).Insert()

$ws3.Range("A3").Value = "$oldGuid.md"
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "$oldGuid.$oldHash.de-de.xlf"
$ws3.Range("E3").Value = "2016-03-22 00:33:37"
$ws3.Range("H3").Value = "0001-01-01 00:00:00"
$ws3.Range("J3").Value = "Include"

$ws3.Range("A2").Value = "$newGuid.md"
$ws3.Range("B2").Value = ".md"
$ws3.Range("C2").Value = "Ready for handoff"
$ws3.Range("D2").Value = "$newGuid.$newHash.de-de.xlf"
$ws3.Range("E2").Value = "2016-03-22 00:34:05"
$ws3.Range("H2").Value = "0001-01-01 00:00:00"
$ws3.Range("J2").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$newGuid.md", "", "", "$newGuid.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$newGuid.$newHash.de-de.xlf", "", "", "$newGuid.$newHash.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/$mdBlobSha/e2e/$oldGuid.md", "", "", "$oldGuid.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deHandoffSha/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldGuid.$oldHash.de-de.xlf", "", "", "$oldGuid.$oldHash.de-de.xlf")

Write-Host "Report for handoff generated."
